$wb = $excel.ActiveWorkbook

# --- Rename "Raffle" sheet to "Events" ---
$wsEvents = $wb.Worksheets.Item("Raffle")
$wsEvents.Name = "Events"

# --- Donation sheet: insert new "Type" column after PaymentType (col F), before CreatedAt ---
$wsDonation = $wb.Worksheets.Item("Donation")
$wsDonation.Columns.Item(7).Insert()
$wsDonation.Range("G1").Value = "Type"
$wsDonation.Range("G2").Value = "Type 1"
$wsDonation.Range("G5").Value = "Type 2"
$wsDonation.Columns.Item(7).ColumnWidth = 11.67

# --- Update selections / active sheet / view state to match the saved session ---

# Donor sheet: selection moves to A10, no longer the visible/topLeft-scrolled tab
$wsDonor = $wb.Worksheets.Item("Donor")
$null = $wsDonor.Range("A10").Select()

# Donation sheet: selection moves to G6 (the new Type column)
$null = $wsDonation.Range("G6").Select()

# Events sheet: selection moves to F15
$null = $wsEvents.Range("F15").Select()

# PhoneType sheet: becomes the active/selected tab, selection at I9
$wsPhoneType = $wb.Worksheets.Item("PhoneType")
$null = $wsPhoneType.Activate()
$null = $wsPhoneType.Range("I9").Select()
